# <tab> moves through textboxes only, not drop downs
#
# The task "tabbing from text should move to next text, not the drop downs"
# (Id 6) has been completed. Move it out of the Active sheet's todo list and
# into the Inactive sheet's Done list (as the newest / top Done entry).

$wb = $excel.ActiveWorkbook
$active = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("Inactive")

# Insert a new row at the top of the Inactive sheet's data (row 2, right
# below the header) for the task that just got marked Done.
$inactive.Rows.Item(2).Insert()

$inactive.Cells.Item(2, 1).Value = 6
$inactive.Cells.Item(2, 2).Value = "tabbing from text should move to next text, not the drop downs"
$inactive.Cells.Item(2, 3).Value = "Done"
$inactive.Cells.Item(2, 4).Value = "Feature"

# Created / Done columns hold dates stored as plain text (e.g. "12/1/2017")
# elsewhere in this workbook, so force text formatting before assigning the
# value - otherwise Excel would silently convert the literal into a real
# date serial number.
$inactive.Cells.Item(2, 5).NumberFormat = "@"
$inactive.Cells.Item(2, 5).Value = "12/1/2017"
$inactive.Cells.Item(2, 6).NumberFormat = "@"
$inactive.Cells.Item(2, 6).Value = "3/3/2018"

# The new row inherited the inserted-above row's (the header's) bold style;
# reset it back to the plain/default style used by the rest of the data rows.
$inactive.Range("A2:F2").Style = "Normal"

# Remove the now-completed task (row 3, Id 6) from the Active sheet; the
# rows below it shift up to fill the gap.
$active.Rows.Item(3).Delete()
